# Apply cryptos list update (Sun Jul 21 12:57:50 UTC 2024, GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.892.11"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "3.498.77"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.13"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.53"
$ws.Range("E6").Value = "  +2.49%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("E8").Value = "  -1.53%  "
$ws.Range("E9").Value = "  +4.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.14"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "4.103.58"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.27"
$ws.Range("E14").Value = "  +4.48%  "
$ws.Range("D15").Value = "66.897.15"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "3.466.21"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.26"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.27"
$ws.Range("E19").Value = "  +2.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "392.61"
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.93"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.17"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("E28").Value = "  -0.26%  "
$ws.Range("E29").Value = "  -3.01%  "
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.05"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "164.29"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.90"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.84"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.25"
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.62"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.840.08"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0734"
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.94"
$ws.Range("E43").Value = "  -1.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.54"
$ws.Range("E44").Value = "  -0.42%  "
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "340.57"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.67"
$ws.Range("E48").Value = "  +3.49%  "
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.841"
$ws.Range("E50").Value = "  -2.21%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.42"
$ws.Range("E51").Value = "  -0.91%  "
